$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.561175525188446
$ws.Range("B1").Value = 1.386369347572327
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.6533123254776
$ws.Range("E1").Value = 1.471076607704163
